$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows 650-652 use the same date style (NumberFormat) as existing column D cells
$dateFormat = $ws.Range("D624").NumberFormat()
$ws.Range("D650").NumberFormat = $dateFormat
$ws.Range("D651").NumberFormat = $dateFormat
$ws.Range("D652").NumberFormat = $dateFormat

# Row 624
$ws.Range("A624").Value = 6
$ws.Range("B624").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C624").Value = "Metropolitana"
$ws.Range("D624").Value = 45075
$ws.Range("E624").Value = 13
$ws.Range("F624").Value = "Fruta"
$ws.Range("G624").Value = 100107
$ws.Range("H624").Value = "Otros"
$ws.Range("I624").Value = 100107011
$ws.Range("J624").Value = "Tuna"
$ws.Range("K624").Value = "Sin especificar"
$ws.Range("L624").Value = "Especial"
$ws.Range("M624").Value = 100
$ws.Range("N624").Value = 24000
$ws.Range("O624").Value = 24000
$ws.Range("P624").Value = 24000
$ws.Range("Q624").Value = "`$/caja 18 kilos"
$ws.Range("R624").Value = "Provincia de Melipilla"
$ws.Range("S624").Value = 1333
$ws.Range("T624").Value = 18

# Row 625
$ws.Range("A625").Value = 6
$ws.Range("B625").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C625").Value = "Metropolitana"
$ws.Range("D625").Value = 45075
$ws.Range("E625").Value = 13
$ws.Range("F625").Value = "Fruta"
$ws.Range("G625").Value = 100107
$ws.Range("H625").Value = "Otros"
$ws.Range("I625").Value = 100107011
$ws.Range("J625").Value = "Tuna"
$ws.Range("K625").Value = "Sin especificar"
$ws.Range("L625").Value = "Primera"
$ws.Range("M625").Value = 80
$ws.Range("N625").Value = 20000
$ws.Range("O625").Value = 20000
$ws.Range("P625").Value = 20000
$ws.Range("Q625").Value = "`$/caja 18 kilos"
$ws.Range("R625").Value = "Provincia de Melipilla"
$ws.Range("S625").Value = 1111
$ws.Range("T625").Value = 18

# Row 626
$ws.Range("A626").Value = 6
$ws.Range("B626").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C626").Value = "Metropolitana"
$ws.Range("D626").Value = 45075
$ws.Range("E626").Value = 13
$ws.Range("F626").Value = "Fruta"
$ws.Range("G626").Value = 100107
$ws.Range("H626").Value = "Otros"
$ws.Range("I626").Value = 100107011
$ws.Range("J626").Value = "Tuna"
$ws.Range("K626").Value = "Sin especificar"
$ws.Range("L626").Value = "Segunda"
$ws.Range("M626").Value = 75
$ws.Range("N626").Value = 15000
$ws.Range("O626").Value = 15000
$ws.Range("P626").Value = 15000
$ws.Range("Q626").Value = "`$/caja 18 kilos"
$ws.Range("R626").Value = "Provincia de Melipilla"
$ws.Range("S626").Value = 833
$ws.Range("T626").Value = 18

# Row 627
$ws.Range("A627").Value = 6
$ws.Range("B627").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C627").Value = "Metropolitana"
$ws.Range("D627").Value = 44722
$ws.Range("E627").Value = 13
$ws.Range("F627").Value = "Fruta"
$ws.Range("G627").Value = 100107
$ws.Range("H627").Value = "Otros"
$ws.Range("I627").Value = 100107011
$ws.Range("J627").Value = "Tuna"
$ws.Range("K627").Value = "Sin especificar"
$ws.Range("L627").Value = "Especial"
$ws.Range("M627").Value = 50
$ws.Range("N627").Value = 22000
$ws.Range("O627").Value = 22000
$ws.Range("P627").Value = 22000
$ws.Range("Q627").Value = "`$/caja 18 kilos"
$ws.Range("R627").Value = "Provincia de Melipilla"
$ws.Range("S627").Value = 1222
$ws.Range("T627").Value = 18

# Row 628
$ws.Range("A628").Value = 6
$ws.Range("B628").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C628").Value = "Metropolitana"
$ws.Range("D628").Value = 44722
$ws.Range("E628").Value = 13
$ws.Range("F628").Value = "Fruta"
$ws.Range("G628").Value = 100107
$ws.Range("H628").Value = "Otros"
$ws.Range("I628").Value = 100107011
$ws.Range("J628").Value = "Tuna"
$ws.Range("K628").Value = "Sin especificar"
$ws.Range("L628").Value = "Primera"
$ws.Range("M628").Value = 75
$ws.Range("N628").Value = 18000
$ws.Range("O628").Value = 18000
$ws.Range("P628").Value = 18000
$ws.Range("Q628").Value = "`$/caja 18 kilos"
$ws.Range("R628").Value = "Provincia de Melipilla"
$ws.Range("S628").Value = 1000
$ws.Range("T628").Value = 18

# Row 629
$ws.Range("A629").Value = 6
$ws.Range("B629").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C629").Value = "Metropolitana"
$ws.Range("D629").Value = 44624
$ws.Range("E629").Value = 13
$ws.Range("F629").Value = "Fruta"
$ws.Range("G629").Value = 100107
$ws.Range("H629").Value = "Otros"
$ws.Range("I629").Value = 100107011
$ws.Range("J629").Value = "Tuna"
$ws.Range("K629").Value = "Sin especificar"
$ws.Range("L629").Value = "Especial"
$ws.Range("M629").Value = 150
$ws.Range("N629").Value = 17000
$ws.Range("O629").Value = 17000
$ws.Range("P629").Value = 17000
$ws.Range("Q629").Value = "`$/caja 18 kilos"
$ws.Range("R629").Value = "Región Metropolitana"
$ws.Range("S629").Value = 944
$ws.Range("T629").Value = 18

# Row 630
$ws.Range("A630").Value = 6
$ws.Range("B630").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C630").Value = "Metropolitana"
$ws.Range("D630").Value = 44624
$ws.Range("E630").Value = 13
$ws.Range("F630").Value = "Fruta"
$ws.Range("G630").Value = 100107
$ws.Range("H630").Value = "Otros"
$ws.Range("I630").Value = 100107011
$ws.Range("J630").Value = "Tuna"
$ws.Range("K630").Value = "Sin especificar"
$ws.Range("L630").Value = "Primera"
$ws.Range("M630").Value = 190
$ws.Range("N630").Value = 14000
$ws.Range("O630").Value = 15000
$ws.Range("P630").Value = 14500
$ws.Range("Q630").Value = "`$/caja 18 kilos"
$ws.Range("R630").Value = "Región Metropolitana"
$ws.Range("S630").Value = 806
$ws.Range("T630").Value = 18

# Row 631
$ws.Range("A631").Value = 6
$ws.Range("B631").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C631").Value = "Metropolitana"
$ws.Range("D631").Value = 44624
$ws.Range("E631").Value = 13
$ws.Range("F631").Value = "Fruta"
$ws.Range("G631").Value = 100107
$ws.Range("H631").Value = "Otros"
$ws.Range("I631").Value = 100107011
$ws.Range("J631").Value = "Tuna"
$ws.Range("K631").Value = "Sin especificar"
$ws.Range("L631").Value = "Segunda"
$ws.Range("M631").Value = 100
$ws.Range("N631").Value = 12000
$ws.Range("O631").Value = 12000
$ws.Range("P631").Value = 12000
$ws.Range("Q631").Value = "`$/caja 18 kilos"
$ws.Range("R631").Value = "Región Metropolitana"
$ws.Range("S631").Value = 667
$ws.Range("T631").Value = 18

# Row 632
$ws.Range("A632").Value = 6
$ws.Range("B632").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C632").Value = "Metropolitana"
$ws.Range("D632").Value = 44624
$ws.Range("E632").Value = 13
$ws.Range("F632").Value = "Fruta"
$ws.Range("G632").Value = 100107
$ws.Range("H632").Value = "Otros"
$ws.Range("I632").Value = 100107011
$ws.Range("J632").Value = "Tuna"
$ws.Range("K632").Value = "Sin especificar"
$ws.Range("L632").Value = "Tercera"
$ws.Range("M632").Value = 100
$ws.Range("N632").Value = 10000
$ws.Range("O632").Value = 10000
$ws.Range("P632").Value = 10000
$ws.Range("Q632").Value = "`$/caja 18 kilos"
$ws.Range("R632").Value = "Región Metropolitana"
$ws.Range("S632").Value = 556
$ws.Range("T632").Value = 18

# Row 633
$ws.Range("A633").Value = 6
$ws.Range("B633").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C633").Value = "Metropolitana"
$ws.Range("D633").Value = 45014
$ws.Range("E633").Value = 13
$ws.Range("F633").Value = "Fruta"
$ws.Range("G633").Value = 100107
$ws.Range("H633").Value = "Otros"
$ws.Range("I633").Value = 100107011
$ws.Range("J633").Value = "Tuna"
$ws.Range("K633").Value = "Sin especificar"
$ws.Range("L633").Value = "Especial"
$ws.Range("M633").Value = 120
$ws.Range("N633").Value = 14000
$ws.Range("O633").Value = 14000
$ws.Range("P633").Value = 14000
$ws.Range("Q633").Value = "`$/caja 18 kilos"
$ws.Range("R633").Value = "Provincia de Melipilla"
$ws.Range("S633").Value = 778
$ws.Range("T633").Value = 18

# Row 634
$ws.Range("A634").Value = 6
$ws.Range("B634").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C634").Value = "Metropolitana"
$ws.Range("D634").Value = 45014
$ws.Range("E634").Value = 13
$ws.Range("F634").Value = "Fruta"
$ws.Range("G634").Value = 100107
$ws.Range("H634").Value = "Otros"
$ws.Range("I634").Value = 100107011
$ws.Range("J634").Value = "Tuna"
$ws.Range("K634").Value = "Sin especificar"
$ws.Range("L634").Value = "Extra (doble especial)"
$ws.Range("M634").Value = 170
$ws.Range("N634").Value = 16000
$ws.Range("O634").Value = 16000
$ws.Range("P634").Value = 16000
$ws.Range("Q634").Value = "`$/caja 18 kilos"
$ws.Range("R634").Value = "Provincia de Melipilla"
$ws.Range("S634").Value = 889
$ws.Range("T634").Value = 18

# Row 635
$ws.Range("A635").Value = 6
$ws.Range("B635").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C635").Value = "Metropolitana"
$ws.Range("D635").Value = 45014
$ws.Range("E635").Value = 13
$ws.Range("F635").Value = "Fruta"
$ws.Range("G635").Value = 100107
$ws.Range("H635").Value = "Otros"
$ws.Range("I635").Value = 100107011
$ws.Range("J635").Value = "Tuna"
$ws.Range("K635").Value = "Sin especificar"
$ws.Range("L635").Value = "Primera"
$ws.Range("M635").Value = 150
$ws.Range("N635").Value = 12000
$ws.Range("O635").Value = 12000
$ws.Range("P635").Value = 12000
$ws.Range("Q635").Value = "`$/caja 18 kilos"
$ws.Range("R635").Value = "Provincia de Melipilla"
$ws.Range("S635").Value = 667
$ws.Range("T635").Value = 18

# Row 636
$ws.Range("A636").Value = 6
$ws.Range("B636").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C636").Value = "Metropolitana"
$ws.Range("D636").Value = 45014
$ws.Range("E636").Value = 13
$ws.Range("F636").Value = "Fruta"
$ws.Range("G636").Value = 100107
$ws.Range("H636").Value = "Otros"
$ws.Range("I636").Value = 100107011
$ws.Range("J636").Value = "Tuna"
$ws.Range("K636").Value = "Sin especificar"
$ws.Range("L636").Value = "Segunda"
$ws.Range("M636").Value = 120
$ws.Range("N636").Value = 10000
$ws.Range("O636").Value = 10000
$ws.Range("P636").Value = 10000
$ws.Range("Q636").Value = "`$/caja 18 kilos"
$ws.Range("R636").Value = "Provincia de Melipilla"
$ws.Range("S636").Value = 556
$ws.Range("T636").Value = 18

# Row 637
$ws.Range("A637").Value = 6
$ws.Range("B637").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C637").Value = "Metropolitana"
$ws.Range("D637").Value = 45014
$ws.Range("E637").Value = 13
$ws.Range("F637").Value = "Fruta"
$ws.Range("G637").Value = 100107
$ws.Range("H637").Value = "Otros"
$ws.Range("I637").Value = 100107011
$ws.Range("J637").Value = "Tuna"
$ws.Range("K637").Value = "Sin especificar"
$ws.Range("L637").Value = "Tercera"
$ws.Range("M637").Value = 100
$ws.Range("N637").Value = 8000
$ws.Range("O637").Value = 8000
$ws.Range("P637").Value = 8000
$ws.Range("Q637").Value = "`$/caja 18 kilos"
$ws.Range("R637").Value = "Provincia de Melipilla"
$ws.Range("S637").Value = 444
$ws.Range("T637").Value = 18

# Row 638
$ws.Range("A638").Value = 6
$ws.Range("B638").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C638").Value = "Metropolitana"
$ws.Range("D638").Value = 44648
$ws.Range("E638").Value = 13
$ws.Range("F638").Value = "Fruta"
$ws.Range("G638").Value = 100107
$ws.Range("H638").Value = "Otros"
$ws.Range("I638").Value = 100107011
$ws.Range("J638").Value = "Tuna"
$ws.Range("K638").Value = "Sin especificar"
$ws.Range("L638").Value = "Especial"
$ws.Range("M638").Value = 200
$ws.Range("N638").Value = 15000
$ws.Range("O638").Value = 15000
$ws.Range("P638").Value = 15000
$ws.Range("Q638").Value = "`$/caja 18 kilos"
$ws.Range("R638").Value = "Región Metropolitana"
$ws.Range("S638").Value = 833
$ws.Range("T638").Value = 18

# Row 639
$ws.Range("A639").Value = 6
$ws.Range("B639").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C639").Value = "Metropolitana"
$ws.Range("D639").Value = 44648
$ws.Range("E639").Value = 13
$ws.Range("F639").Value = "Fruta"
$ws.Range("G639").Value = 100107
$ws.Range("H639").Value = "Otros"
$ws.Range("I639").Value = 100107011
$ws.Range("J639").Value = "Tuna"
$ws.Range("K639").Value = "Sin especificar"
$ws.Range("L639").Value = "Primera"
$ws.Range("M639").Value = 200
$ws.Range("N639").Value = 12000
$ws.Range("O639").Value = 12000
$ws.Range("P639").Value = 12000
$ws.Range("Q639").Value = "`$/caja 18 kilos"
$ws.Range("R639").Value = "Región Metropolitana"
$ws.Range("S639").Value = 667
$ws.Range("T639").Value = 18

# Row 640
$ws.Range("A640").Value = 6
$ws.Range("B640").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C640").Value = "Metropolitana"
$ws.Range("D640").Value = 44648
$ws.Range("E640").Value = 13
$ws.Range("F640").Value = "Fruta"
$ws.Range("G640").Value = 100107
$ws.Range("H640").Value = "Otros"
$ws.Range("I640").Value = 100107011
$ws.Range("J640").Value = "Tuna"
$ws.Range("K640").Value = "Sin especificar"
$ws.Range("L640").Value = "Segunda"
$ws.Range("M640").Value = 200
$ws.Range("N640").Value = 9000
$ws.Range("O640").Value = 9000
$ws.Range("P640").Value = 9000
$ws.Range("Q640").Value = "`$/caja 18 kilos"
$ws.Range("R640").Value = "Región Metropolitana"
$ws.Range("S640").Value = 500
$ws.Range("T640").Value = 18

# Row 641
$ws.Range("A641").Value = 6
$ws.Range("B641").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C641").Value = "Metropolitana"
$ws.Range("D641").Value = 44651
$ws.Range("E641").Value = 13
$ws.Range("F641").Value = "Fruta"
$ws.Range("G641").Value = 100107
$ws.Range("H641").Value = "Otros"
$ws.Range("I641").Value = 100107011
$ws.Range("J641").Value = "Tuna"
$ws.Range("K641").Value = "Sin especificar"
$ws.Range("L641").Value = "Especial"
$ws.Range("M641").Value = 250
$ws.Range("N641").Value = 13000
$ws.Range("O641").Value = 13000
$ws.Range("P641").Value = 13000
$ws.Range("Q641").Value = "`$/caja 18 kilos"
$ws.Range("R641").Value = "Provincia de Limarí"
$ws.Range("S641").Value = 722
$ws.Range("T641").Value = 18

# Row 642
$ws.Range("A642").Value = 6
$ws.Range("B642").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C642").Value = "Metropolitana"
$ws.Range("D642").Value = 44651
$ws.Range("E642").Value = 13
$ws.Range("F642").Value = "Fruta"
$ws.Range("G642").Value = 100107
$ws.Range("H642").Value = "Otros"
$ws.Range("I642").Value = 100107011
$ws.Range("J642").Value = "Tuna"
$ws.Range("K642").Value = "Sin especificar"
$ws.Range("L642").Value = "Especial"
$ws.Range("M642").Value = 250
$ws.Range("N642").Value = 13000
$ws.Range("O642").Value = 13000
$ws.Range("P642").Value = 13000
$ws.Range("Q642").Value = "`$/caja 18 kilos"
$ws.Range("R642").Value = "Región Metropolitana"
$ws.Range("S642").Value = 722
$ws.Range("T642").Value = 18

# Row 643
$ws.Range("A643").Value = 6
$ws.Range("B643").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C643").Value = "Metropolitana"
$ws.Range("D643").Value = 44651
$ws.Range("E643").Value = 13
$ws.Range("F643").Value = "Fruta"
$ws.Range("G643").Value = 100107
$ws.Range("H643").Value = "Otros"
$ws.Range("I643").Value = 100107011
$ws.Range("J643").Value = "Tuna"
$ws.Range("K643").Value = "Sin especificar"
$ws.Range("L643").Value = "Primera"
$ws.Range("M643").Value = 300
$ws.Range("N643").Value = 10000
$ws.Range("O643").Value = 10000
$ws.Range("P643").Value = 10000
$ws.Range("Q643").Value = "`$/caja 18 kilos"
$ws.Range("R643").Value = "Provincia de Limarí"
$ws.Range("S643").Value = 556
$ws.Range("T643").Value = 18

# Row 644
$ws.Range("A644").Value = 6
$ws.Range("B644").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C644").Value = "Metropolitana"
$ws.Range("D644").Value = 44651
$ws.Range("E644").Value = 13
$ws.Range("F644").Value = "Fruta"
$ws.Range("G644").Value = 100107
$ws.Range("H644").Value = "Otros"
$ws.Range("I644").Value = 100107011
$ws.Range("J644").Value = "Tuna"
$ws.Range("K644").Value = "Sin especificar"
$ws.Range("L644").Value = "Primera"
$ws.Range("M644").Value = 325
$ws.Range("N644").Value = 10000
$ws.Range("O644").Value = 10000
$ws.Range("P644").Value = 10000
$ws.Range("Q644").Value = "`$/caja 18 kilos"
$ws.Range("R644").Value = "Región Metropolitana"
$ws.Range("S644").Value = 556
$ws.Range("T644").Value = 18

# Row 645
$ws.Range("A645").Value = 6
$ws.Range("B645").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C645").Value = "Metropolitana"
$ws.Range("D645").Value = 44651
$ws.Range("E645").Value = 13
$ws.Range("F645").Value = "Fruta"
$ws.Range("G645").Value = 100107
$ws.Range("H645").Value = "Otros"
$ws.Range("I645").Value = 100107011
$ws.Range("J645").Value = "Tuna"
$ws.Range("K645").Value = "Sin especificar"
$ws.Range("L645").Value = "Segunda"
$ws.Range("M645").Value = 270
$ws.Range("N645").Value = 8000
$ws.Range("O645").Value = 8000
$ws.Range("P645").Value = 8000
$ws.Range("Q645").Value = "`$/caja 18 kilos"
$ws.Range("R645").Value = "Provincia de Limarí"
$ws.Range("S645").Value = 444
$ws.Range("T645").Value = 18

# Row 646
$ws.Range("A646").Value = 6
$ws.Range("B646").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C646").Value = "Metropolitana"
$ws.Range("D646").Value = 44651
$ws.Range("E646").Value = 13
$ws.Range("F646").Value = "Fruta"
$ws.Range("G646").Value = 100107
$ws.Range("H646").Value = "Otros"
$ws.Range("I646").Value = 100107011
$ws.Range("J646").Value = "Tuna"
$ws.Range("K646").Value = "Sin especificar"
$ws.Range("L646").Value = "Segunda"
$ws.Range("M646").Value = 270
$ws.Range("N646").Value = 8000
$ws.Range("O646").Value = 8000
$ws.Range("P646").Value = 8000
$ws.Range("Q646").Value = "`$/caja 18 kilos"
$ws.Range("R646").Value = "Región Metropolitana"
$ws.Range("S646").Value = 444
$ws.Range("T646").Value = 18

# Row 647
$ws.Range("A647").Value = 6
$ws.Range("B647").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C647").Value = "Metropolitana"
$ws.Range("D647").Value = 44988
$ws.Range("E647").Value = 13
$ws.Range("F647").Value = "Fruta"
$ws.Range("G647").Value = 100107
$ws.Range("H647").Value = "Otros"
$ws.Range("I647").Value = 100107011
$ws.Range("J647").Value = "Tuna"
$ws.Range("K647").Value = "Sin especificar"
$ws.Range("L647").Value = "Especial"
$ws.Range("M647").Value = 275
$ws.Range("N647").Value = 17000
$ws.Range("O647").Value = 17000
$ws.Range("P647").Value = 17000
$ws.Range("Q647").Value = "`$/caja 18 kilos"
$ws.Range("R647").Value = "Región Metropolitana"
$ws.Range("S647").Value = 944
$ws.Range("T647").Value = 18

# Row 648
$ws.Range("A648").Value = 6
$ws.Range("B648").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C648").Value = "Metropolitana"
$ws.Range("D648").Value = 44988
$ws.Range("E648").Value = 13
$ws.Range("F648").Value = "Fruta"
$ws.Range("G648").Value = 100107
$ws.Range("H648").Value = "Otros"
$ws.Range("I648").Value = 100107011
$ws.Range("J648").Value = "Tuna"
$ws.Range("K648").Value = "Sin especificar"
$ws.Range("L648").Value = "Primera"
$ws.Range("M648").Value = 275
$ws.Range("N648").Value = 13000
$ws.Range("O648").Value = 13000
$ws.Range("P648").Value = 13000
$ws.Range("Q648").Value = "`$/caja 18 kilos"
$ws.Range("R648").Value = "Región Metropolitana"
$ws.Range("S648").Value = 722
$ws.Range("T648").Value = 18

# Row 649
$ws.Range("A649").Value = 6
$ws.Range("B649").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C649").Value = "Metropolitana"
$ws.Range("D649").Value = 44988
$ws.Range("E649").Value = 13
$ws.Range("F649").Value = "Fruta"
$ws.Range("G649").Value = 100107
$ws.Range("H649").Value = "Otros"
$ws.Range("I649").Value = 100107011
$ws.Range("J649").Value = "Tuna"
$ws.Range("K649").Value = "Sin especificar"
$ws.Range("L649").Value = "Segunda"
$ws.Range("M649").Value = 275
$ws.Range("N649").Value = 10000
$ws.Range("O649").Value = 10000
$ws.Range("P649").Value = 10000
$ws.Range("Q649").Value = "`$/caja 18 kilos"
$ws.Range("R649").Value = "Región Metropolitana"
$ws.Range("S649").Value = 556
$ws.Range("T649").Value = 18

# Row 650
$ws.Range("A650").Value = 6
$ws.Range("B650").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C650").Value = "Metropolitana"
$ws.Range("D650").Value = 44999
$ws.Range("E650").Value = 13
$ws.Range("F650").Value = "Fruta"
$ws.Range("G650").Value = 100107
$ws.Range("H650").Value = "Otros"
$ws.Range("I650").Value = 100107011
$ws.Range("J650").Value = "Tuna"
$ws.Range("K650").Value = "Sin especificar"
$ws.Range("L650").Value = "Especial"
$ws.Range("M650").Value = 275
$ws.Range("N650").Value = 14000
$ws.Range("O650").Value = 14000
$ws.Range("P650").Value = 14000
$ws.Range("Q650").Value = "`$/caja 18 kilos"
$ws.Range("R650").Value = "Provincia de Melipilla"
$ws.Range("S650").Value = 778
$ws.Range("T650").Value = 18

# Row 651
$ws.Range("A651").Value = 6
$ws.Range("B651").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C651").Value = "Metropolitana"
$ws.Range("D651").Value = 44999
$ws.Range("E651").Value = 13
$ws.Range("F651").Value = "Fruta"
$ws.Range("G651").Value = 100107
$ws.Range("H651").Value = "Otros"
$ws.Range("I651").Value = 100107011
$ws.Range("J651").Value = "Tuna"
$ws.Range("K651").Value = "Sin especificar"
$ws.Range("L651").Value = "Primera"
$ws.Range("M651").Value = 275
$ws.Range("N651").Value = 12000
$ws.Range("O651").Value = 12000
$ws.Range("P651").Value = 12000
$ws.Range("Q651").Value = "`$/caja 18 kilos"
$ws.Range("R651").Value = "Provincia de Melipilla"
$ws.Range("S651").Value = 667
$ws.Range("T651").Value = 18

# Row 652
$ws.Range("A652").Value = 6
$ws.Range("B652").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C652").Value = "Metropolitana"
$ws.Range("D652").Value = 44999
$ws.Range("E652").Value = 13
$ws.Range("F652").Value = "Fruta"
$ws.Range("G652").Value = 100107
$ws.Range("H652").Value = "Otros"
$ws.Range("I652").Value = 100107011
$ws.Range("J652").Value = "Tuna"
$ws.Range("K652").Value = "Sin especificar"
$ws.Range("L652").Value = "Segunda"
$ws.Range("M652").Value = 275
$ws.Range("N652").Value = 10000
$ws.Range("O652").Value = 10000
$ws.Range("P652").Value = 10000
$ws.Range("Q652").Value = "`$/caja 18 kilos"
$ws.Range("R652").Value = "Provincia de Melipilla"
$ws.Range("S652").Value = 556
$ws.Range("T652").Value = 18
